$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 3: add P3 as an empty bottom-bordered cell (same format as A3/C3/O3) ---
$ws.Range("O3").Copy() | Out-Null
$ws.Range("P3").PasteSpecial($xlPasteFormats) | Out-Null

# --- Row 4: header year 2021, same format as O4 ---
$ws.Range("O4").Copy() | Out-Null
$ws.Range("P4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("P4").Value = 2021

# --- Row 5: total victims for 2021 ---
$ws.Range("O8").Copy() | Out-Null
$ws.Range("P5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("P5").Value = 9038

# --- Row 6: section header "By sex:" row, blank numeric cell ---
$ws.Range("O6").Copy() | Out-Null
$ws.Range("P6").PasteSpecial($xlPasteFormats) | Out-Null

# --- Row 7: Women ---
$ws.Range("O6").Copy() | Out-Null
$ws.Range("P7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("P7").Value = 8587

# --- Row 8: Men ---
$ws.Range("O6").Copy() | Out-Null
$ws.Range("P8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("P8").Value = 451

# --- Row 9: section header "By level of education:" row, blank numeric cell ---
$ws.Range("O6").Copy() | Out-Null
$ws.Range("P9").PasteSpecial($xlPasteFormats) | Out-Null

# --- Rows 10-24: data rows with no 2021 figures available ("…") ---
$ws.Range("O10").Copy() | Out-Null
for ($r = 10; $r -le 24; $r++) {
    $ws.Cells.Item($r, 16).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Cells.Item($r, 16).Value = "…"
}

# --- Row 25: last data row, bottom-bordered ("…") ---
$ws.Range("O25").Copy() | Out-Null
$ws.Range("P25").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("P25").Value = "…"

$excel.CutCopyMode = $false

# Match the recorded selection left behind in the saved sheet view
$ws.Range("Q4").Select() | Out-Null
